$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new table (row 13-14, merged vertically)
$ws.Range("A13").Value = "Array"
$ws.Range("A14").Value = "size"
$ws.Range("B13").Value = "Bubble Sort"
$ws.Range("C13").Value = "Insertion Sort"
$ws.Range("D13").Value = "Selection Sort"
$ws.Range("E13").Value = "Quick S"
$ws.Range("F13").Value = "Merge S"
$ws.Range("G13").Value = "Heap S."

# Data rows 15-22
$data = @(
    @(50, 0, 0, 0, 0, 0, 0),
    @(100, 1, 0, 0, 0, 0, 0),
    @(500, 6.3, 3.2, 3, 1, 0, 1),
    @(1000, 32.13, 13.52, 13, 1.01, 1.02, 1.01),
    @(5000, 778.71, 348.76, 323.16, 7.1, 7, 11.14),
    @(10000, 3249.78, 1319.78, 1252.38, 12.27, 14.61, 21.61),
    @(50000, 82193.01, 33936.97, 33050.45, 72.11, 84.18, 130.78),
    @(100000, 338244.59, 150274.49, 137543.03, 173.95, 172.97, 294.93)
)

$r = 15
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
